$d = $word.ActiveDocument

# --- Step 1: Append new runs to the end of the "Usage Forecast" body paragraph (paragraph 20) ---
$p20 = $d.Paragraphs.Item(20)
$r20 = $p20.Range
$r20.InsertAfter(' The high point should be the learnability, satisfaction, and memorability')
$r20 = $p20.Range
$r20.InsertAfter('.')
$r20 = $p20.Range
$r20.InsertAfter(' The achievement system should lead the user to learning how to use this system')
$r20 = $p20.Range
$r20.InsertAfter('. The achievements will be made in a step by step manner to get the user to learn it at their own pace. The achievement they unlock will be something that stays with them, even if it was an achievement for which they did the action only once that memory will stay with them and serve as a constant reminder that such action is possible and is something they have done before.')
$r20 = $p20.Range
$r20.InsertAfter(' Even if they didn’t quite understand how they did it the first time they will at least have some sort of idea on how they did it and can look to the community to help then fill in the missing steps they may have forgotten.')
$r20 = $p20.Range
$r20.InsertAfter(' The satisfaction should be one of the highest points because the user will always be as satisfied as they want to be.')
$r20 = $p20.Range
$r20.InsertAfter(' With all of the customization options and the ')
$r20 = $p20.Range
$r20.InsertAfter('achievements')
$r20 = $p20.Range
$r20.InsertAfter(' teaching how to customize at will the users should be using their own dream interface web browser.')
$r20 = $p20.Range
$r20.InsertAfter(' ')
$r20 = $p20.Range
$r20.InsertAfter('Overall the browser should be well received by the general public.')

# --- Step 2: Insert a new (non-bold) paragraph after paragraph 20 and fill it with the "In conclusion" text ---
$p20b = $d.Paragraphs.Item(20)
$r20b = $p20b.Range
$r20b.InsertParagraphAfter()
$pInConclusion = $d.Paragraphs.Item(21)
$rInConclusion = $pInConclusion.Range
$rInConclusion.InsertAfter('In conclusion this browser should be seen as a beginner’s web browser ')
$pInConclusion = $d.Paragraphs.Item(21)
$rInConclusion = $pInConclusion.Range
$rInConclusion.InsertAfter('but also a great tool for developers.')
$pInConclusion = $d.Paragraphs.Item(21)
$rInConclusion = $pInConclusion.Range
$rInConclusion.InsertAfter(' It has a full customizability that a beginner can grasp the basics ')
$pInConclusion = $d.Paragraphs.Item(21)
$rInConclusion = $pInConclusion.Range
$rInConclusion.InsertAfter('of and')
$pInConclusion = $d.Paragraphs.Item(21)
$rInConclusion = $pInConclusion.Range
$rInConclusion.InsertAfter(' ')
$pInConclusion = $d.Paragraphs.Item(21)
$rInConclusion = $pInConclusion.Range
$rInConclusion.InsertAfter('a more accustomed user ')
$pInConclusion = $d.Paragraphs.Item(21)
$rInConclusion = $pInConclusion.Range
$rInConclusion.InsertAfter('can edit to a much lower level.')
$pInConclusion = $d.Paragraphs.Item(21)
$rInConclusion = $pInConclusion.Range
$rInConclusion.InsertAfter(' ')
$pInConclusion = $d.Paragraphs.Item(21)
$rInConclusion = $pInConclusion.Range
$rInConclusion.InsertAfter('The brick design is meant to be something familiar to most and ')
$pInConclusion = $d.Paragraphs.Item(21)
$rInConclusion = $pInConclusion.Range
$rInConclusion.InsertAfter('communicates the idea of ')
$pInConclusion = $d.Paragraphs.Item(21)
$rInConclusion = $pInConclusion.Range
$rInConclusion.InsertAfter('having the ability to ')
$pInConclusion = $d.Paragraphs.Item(21)
$rInConclusion = $pInConclusion.Range
$rInConclusion.InsertAfter('be p')
$pInConclusion = $d.Paragraphs.Item(21)
$rInConclusion = $pInConclusion.Range
$rInConclusion.InsertAfter('ut together.')
$pInConclusion = $d.Paragraphs.Item(21)
$rInConclusion = $pInConclusion.Range
$rInConclusion.InsertAfter(' Each brick is meant work on its own and also work with others for a greater brick and different experience.')

# --- Step 3: Insert a new bold "Conclusion" heading paragraph before the "In conclusion" paragraph ---
$pInConclusion2 = $d.Paragraphs.Item(21)
$rInConclusion2 = $pInConclusion2.Range
$rInConclusion2.InsertParagraphBefore()
$pConclusion = $d.Paragraphs.Item(21)
$rConclusion = $pConclusion.Range
$rConclusion.InsertAfter('Conclusion')
$pConclusion2 = $d.Paragraphs.Item(21)
$rConclusion2 = $pConclusion2.Range
$boldRange = $d.Range($rConclusion2.Start, $rConclusion2.Start + 10)
$boldRange.Bold = 1

# --- Step 4: Delete the old trailing paragraph ("Make sure you copy paste...") ---
$pOld = $d.Paragraphs.Item(23)
$rOld = $pOld.Range
$rOld.Delete()

# --- Step 5: Re-create the _GoBack bookmark at the very end of the "In conclusion" paragraph (collapsed) ---
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$rLast = $pLast.Range
$rLast.InsertAfter("Z")
$pLast2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$rLast2 = $pLast2.Range
$bmRange = $d.Range($rLast2.End - 2, $rLast2.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$pLast3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$rLast3 = $pLast3.Range
$zRange = $d.Range($rLast3.End - 2, $rLast3.End - 1)
$zRange.Delete()
